$wb = $excel.ActiveWorkbook

# --- Sheet2: populate the new data table (A1:B7) ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws2.Range("A1").Value = "H1"
$ws2.Range("B1").Value = "H2"

$ws2.Range("A2").Value = "id"
$ws2.Range("B2").Value = 101

$ws2.Range("A3").Value = "name"
$ws2.Range("B3").Value = "john"

$ws2.Range("A4").Value = "tag-id"
$ws2.Range("B4").Value = 11

$ws2.Range("A5").Value = "tag-name"
$ws2.Range("B5").Value = "peter1"

$ws2.Range("A6").Value = "tag-id"
$ws2.Range("B6").Value = 12

$ws2.Range("A7").Value = "tag-name"
$ws2.Range("B7").Value = "peter2"

# Page setup for Sheet2 (paperSize=9 / Letter, portrait orientation)
$ps2 = $ws2.PageSetup
$ps2.PaperSize = 9
$ps2.Orientation = 1

# --- Active sheet / tab selection moves from Sheet1 to Sheet2 ---
$ws2.Activate() | Out-Null
$ws2.Range("B8").Select() | Out-Null

# --- Window geometry (best effort) ---
$win = $wb.Windows.Item(1)
$win.Left = 2229
$win.Top = 2229
$win.Width = 16457
$win.Height = 9548
